# Append a new trade record as row 6 on the (single) worksheet.
# Columns: A=Date, B=Profitable, C=Principle, D=Start Principle,
#          E=BuyPrice, F=SellPrice, G=IsShortSell, H=Price Change %, I=Strong trade

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A6: Date (stored as serial date, formatted like the existing A3:A5 date cells)
$ws.Range("A6").Value = 42647.681817129633
$ws.Range("A6").NumberFormat = "m/d/yy h:mm"

# B6: Profitable (boolean)
$ws.Range("B6").Value = $false

# C6: Principle
$ws.Range("C6").Value = 9971.89

# D6: Start Principle
$ws.Range("D6").Value = 10013.450000000001

# E6: BuyPrice
$ws.Range("E6").Value = 18.12

# F6: SellPrice
$ws.Range("F6").Value = 17.97

# G6: IsShortSell (boolean, same date-number-format as column A/G above, matching existing rows)
$ws.Range("G6").Value = $false
$ws.Range("G6").NumberFormat = "m/d/yy h:mm"

# H6: Price Change %
$ws.Range("H6").Value = -0.83

# I6: Strong trade (boolean)
$ws.Range("I6").Value = $false
